$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 392.83334
$ws.Range("I28").Value = 258.15384
$ws.Range("J28").Value = 552
$ws.Range("K28").Value = 258.15384
$ws.Range("L28").Value = 552
$ws.Range("M28").Value = 226.84616
$ws.Range("N28").Value = -1522

# Row 86
$ws.Range("H86").Value = 5719.7
$ws.Range("I86").Value = 6877.15
$ws.Range("J86").Value = 3404.8
$ws.Range("K86").Value = 6877.15
$ws.Range("L86").Value = 3404.8
$ws.Range("M86").Value = -5754.15
$ws.Range("N86").Value = -5650.8

# Row 89
$ws.Range("H89").Value = 5719.7
$ws.Range("I89").Value = 6877.15
$ws.Range("J89").Value = 3404.8
$ws.Range("K89").Value = 34385.75
$ws.Range("L89").Value = 17024
$ws.Range("M89").Value = -28769.75
$ws.Range("N89").Value = -28256

# Row 98
$ws.Range("H98").Value = 1661
$ws.Range("I98").Value = 1225
$ws.Range("J98").Value = 2184.2
$ws.Range("K98").Value = 1225
$ws.Range("L98").Value = 2184.2
$ws.Range("M98").Value = 273
$ws.Range("N98").Value = -5180.2

# Row 107
$ws.Range("H107").Value = 38690.117
$ws.Range("I107").Value = 43699.78
$ws.Range("J107").Value = 282.66666
$ws.Range("K107").Value = 43699.78
$ws.Range("L107").Value = 282.66666
$ws.Range("M107").Value = -41779.78
$ws.Range("N107").Value = -4122.66666

# Row 122
$ws.Range("H122").Value = 1661
$ws.Range("I122").Value = 1225
$ws.Range("J122").Value = 2184.2
$ws.Range("K122").Value = 3675
$ws.Range("L122").Value = 6552.599999999999
$ws.Range("M122").Value = -1225
$ws.Range("N122").Value = -11452.6

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21823.906
$ws.Range("I32").Value = 4241.0654
$ws.Range("J32").Value = 104328
$ws.Range("K32").Value = 4241.0654
$ws.Range("L32").Value = 104328
$ws.Range("M32").Value = -3954.0654
$ws.Range("N32").Value = -104902

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2485.2144
$ws.Range("I86").Value = 2862
$ws.Range("J86").Value = 1982.8334
$ws.Range("K86").Value = 2862
$ws.Range("L86").Value = 1982.8334
$ws.Range("M86").Value = -1739
$ws.Range("N86").Value = -4228.8334

# Row 89
$ws.Range("H89").Value = 2485.2144
$ws.Range("I89").Value = 2862
$ws.Range("J89").Value = 1982.8334
$ws.Range("K89").Value = 14310
$ws.Range("L89").Value = 9914.166999999999
$ws.Range("M89").Value = -8694
$ws.Range("N89").Value = -21146.167

# Row 134
$ws.Range("H134").Value = 2091.75
$ws.Range("I134").Value = 2060.3462
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 6181.0386
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -3646.0386
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2451.016
$ws.Range("I31").Value = 1946.84
$ws.Range("K31").Value = 1946.84
$ws.Range("M31").Value = -1651.84

# Row 34
$ws.Range("H34").Value = 2451.016
$ws.Range("I34").Value = 1946.84
$ws.Range("K34").Value = 1946.84
$ws.Range("M34").Value = -1744.84

# Row 105
$ws.Range("H105").Value = 5413.25
$ws.Range("I105").Value = 5800.143
$ws.Range("J105").Value = 5112.3335
$ws.Range("K105").Value = 5800.143
$ws.Range("L105").Value = 5112.3335
$ws.Range("M105").Value = -4053.143
$ws.Range("N105").Value = -8606.333500000001

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1990.3654
$ws.Range("I68").Value = 1827.9131
$ws.Range("J68").Value = 2119.2068
$ws.Range("K68").Value = 5483.7393
$ws.Range("L68").Value = 6357.6204
$ws.Range("M68").Value = -4672.7393
$ws.Range("N68").Value = -7979.6204

# Row 71
$ws.Range("H71").Value = 1990.3654
$ws.Range("I71").Value = 1827.9131
$ws.Range("J71").Value = 2119.2068
$ws.Range("K71").Value = 16451.2179
$ws.Range("L71").Value = 19072.8612
$ws.Range("M71").Value = -12395.2179
$ws.Range("N71").Value = -27184.8612

# Row 101
$ws.Range("H101").Value = 18750
$ws.Range("J101").Value = 18750
$ws.Range("L101").Value = 56250
$ws.Range("N101").Value = -61118

# Row 105
$ws.Range("H105").Value = 7460.3335
$ws.Range("I105").Value = 1950
$ws.Range("J105").Value = 8149.125
$ws.Range("K105").Value = 5850
$ws.Range("L105").Value = 24447.375
$ws.Range("M105").Value = -3229
$ws.Range("N105").Value = -29689.375

# Row 107
$ws.Range("H107").Value = 864.4286
$ws.Range("I107").Value = 498.21622
$ws.Range("J107").Value = 1577.579
$ws.Range("K107").Value = 1494.64866
$ws.Range("L107").Value = 4732.737
$ws.Range("M107").Value = 425.3513399999999
$ws.Range("N107").Value = -8572.737000000001

# Row 110
$ws.Range("H110").Value = 7000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 21000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -29180

# Row 131
$ws.Range("H131").Value = 899.75
$ws.Range("I131").Value = 455.7143
$ws.Range("J131").Value = 1138.8462
$ws.Range("K131").Value = 1367.1429
$ws.Range("L131").Value = 3416.5386
$ws.Range("M131").Value = 3672.8571
$ws.Range("N131").Value = -13496.5386

# Row 137
$ws.Range("H137").Value = 8581.538
$ws.Range("I137").Value = 2265
$ws.Range("J137").Value = 11388.889
$ws.Range("K137").Value = 6795
$ws.Range("L137").Value = 34166.667
$ws.Range("M137").Value = -1695
$ws.Range("N137").Value = -44366.667

$ws = $wb.Worksheets.Item("GSM")
# Row 133
$ws.Range("H133").Value = 38600
$ws.Range("J133").Value = 38600
$ws.Range("L133").Value = 38600
$ws.Range("N133").Value = -48720

$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 37607.5
$ws.Range("J92").Value = 37607.5
$ws.Range("L92").Value = 37607.5
$ws.Range("N92").Value = -42599.5

# Row 122
$ws.Range("H122").Value = 3140
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 5990
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 17970
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -22870

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1371.5193
$ws.Range("I136").Value = 1181.6428
$ws.Range("J136").Value = 2169
$ws.Range("K136").Value = 3544.9284
$ws.Range("L136").Value = 6507
$ws.Range("M136").Value = -994.9284000000002
$ws.Range("N136").Value = -11607
